# feat: add 2022-Q3 data
#
# - "总计" (summary) sheet: insert a new leading data row for 2022-Q3
#   (持有数量=5, 持有市值=0.05) and push the existing 2021-Q3 / 2021-Q2
#   rows down by one.
# - Add a new "2022-Q3" worksheet (cloned from "2021-Q3" so it keeps the
#   same header styling / borders) placed right after "总计", populated
#   with the 2022-Q3 fund holdings table.

$wb = $excel.ActiveWorkbook

function Set-TextValue($rng, $val) {
    # Force a value to be stored as TEXT (even when it looks numeric, e.g.
    # a leading-zero fund code or a "29.80" figure that must keep its
    # trailing zero) without leaving a residual NumberFormat/style behind.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# ---------------------------------------------------------------------
# 1. "总计" sheet — shift rows 2:3 down to 3:4, insert the 2022-Q3 row.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Range("A3:D3").Copy($summary.Range("A4:D4"))
$summary.Range("A2:D2").Copy($summary.Range("A3:D3"))

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 5
$summary.Range("D2").Value = 0.05

$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2

# ---------------------------------------------------------------------
# 2. New "2022-Q3" sheet — clone "2021-Q3" (same column layout/styling),
#    rename, reposition right after "总计", then overwrite with the new
#    quarter's data.
# ---------------------------------------------------------------------
$q3_2021 = $wb.Worksheets.Item("2021-Q3")
$q3_2021.Copy($null, $summary)

$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

# the source sheet had 6 data rows (2021-Q3); 2022-Q3 only has 5, so drop
# the now-unused last row.
$newSheet.Rows.Item(7).Delete()

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("A2").Value = 0
Set-TextValue $newSheet.Range("B2") "003842"
$newSheet.Range("C2").Value = "中邮景泰灵活配置混合A"
Set-TextValue $newSheet.Range("D2") "2.09"
Set-TextValue $newSheet.Range("E2") "35.39"
Set-TextValue $newSheet.Range("F2") "1.34"
Set-TextValue $newSheet.Range("G2") "0.0280"
$newSheet.Range("H2").Value = 7

$newSheet.Range("A3").Value = 1
Set-TextValue $newSheet.Range("B3") "004244"
$newSheet.Range("C3").Value = "东方周期优选灵活配置混合"
Set-TextValue $newSheet.Range("D3") "0.35"
Set-TextValue $newSheet.Range("E3") "90.61"
Set-TextValue $newSheet.Range("F3") "4.21"
Set-TextValue $newSheet.Range("G3") "0.0147"
$newSheet.Range("H3").Value = 10

$newSheet.Range("A4").Value = 2
Set-TextValue $newSheet.Range("B4") "159787"
$newSheet.Range("C4").Value = "易方达中证全指建筑材料ETF"
Set-TextValue $newSheet.Range("D4") "0.17"
Set-TextValue $newSheet.Range("E4") "94.24"
Set-TextValue $newSheet.Range("F4") "2.66"
Set-TextValue $newSheet.Range("G4") "0.0045"
$newSheet.Range("H4").Value = 10

$newSheet.Range("A5").Value = 3
Set-TextValue $newSheet.Range("B5") "001430"
$newSheet.Range("C5").Value = "中邮乐享收益灵活配置混合"
Set-TextValue $newSheet.Range("D5") "0.13"
Set-TextValue $newSheet.Range("E5") "29.80"
Set-TextValue $newSheet.Range("F5") "1.12"
Set-TextValue $newSheet.Range("G5") "0.0015"
$newSheet.Range("H5").Value = 10

$newSheet.Range("A6").Value = 4
Set-TextValue $newSheet.Range("B6") "003843"
$newSheet.Range("C6").Value = "中邮景泰灵活配置混合C"
Set-TextValue $newSheet.Range("D6") "0.06"
Set-TextValue $newSheet.Range("E6") "35.39"
Set-TextValue $newSheet.Range("F6") "1.34"
Set-TextValue $newSheet.Range("G6") "0.0008"
$newSheet.Range("H6").Value = 7


# Keep "总计" as the active/selected sheet (untouched by the diff -
# bookViews/activeTab stayed at 0), rather than leaving the newly-created
# "2022-Q3" sheet focused.
$summary.Activate()

Write-Output "2022-Q3 sheet added"
